$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "b" labels in column A (rows 2-7) down by one index:
# b1 -> b0, b2 -> b1, b3 -> b2, b4 -> b3, b5 -> b4, b6 -> b5
for ($i = 2; $i -le 7; $i++) {
    $ws.Cells.Item($i, 1).Value = "b" + ($i - 2)
}
